# Insert a new data row at row 31 (pushes existing rows 31-72 down to 32-73)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 31; this shifts rows 31..72 -> 32..73
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the new record.
$ws.Cells.Item(31, 1).Value  = 3
$ws.Cells.Item(31, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(31, 3).Value  = "Coquimbo"
$ws.Cells.Item(31, 4).Value  = 44894
$ws.Cells.Item(31, 5).Value  = 5
$ws.Cells.Item(31, 6).Value  = 100112022
$ws.Cells.Item(31, 7).Value  = "Arveja Verde"
$ws.Cells.Item(31, 8).Value  = "Perfection"
$ws.Cells.Item(31, 9).Value  = "Primera"
$ws.Cells.Item(31, 10).Value = 73
$ws.Cells.Item(31, 11).Value = 22000
$ws.Cells.Item(31, 12).Value = 23000
$ws.Cells.Item(31, 13).Value = 22521
$ws.Cells.Item(31, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(31, 15).Value = "Región Metropolitana"
$ws.Cells.Item(31, 16).Value = 901
$ws.Cells.Item(31, 17).Value = 25
$ws.Cells.Item(31, 18).Value = "Hortaliza"
